$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the AIM/closing-date data in column J (rows 2-11), which is no
# longer read correctly from the AIM file.
$ws.Range("J2:J11").ClearContents()

# Update the selection to reflect the columns that were just cleared.
$ws.Range("J2:J11").Select()
